$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# New results for the "Refactoring" action
$ws.Range("A29").Value = "Refactoring"
$ws.Range("B29").Value = "Pair Programming (all)"

$ws.Range("A30").Value = "Everything right"
$ws.Range("B30").Value = 17

$ws.Range("A31").Value = "Everything wrong"
$ws.Range("B31").Value = 26

# Scroll the view down and select the next empty row, as in the saved workbook
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A32").Select()
